$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update vehicle 1 (row 2) Patente/Motor/Chasis ---
$ws.Range("H2").Value = "ZZZ518"
$ws.Range("I2").Value = "ABC0987AX314"
$ws.Range("J2").Value = "MMAA09XFGS309"

# --- Update vehicle 2 (row 3) Patente/Motor/Chasis ---
$ws.Range("H3").Value = "ZZZ519"
$ws.Range("I3").Value = "ABC0987AX315"
$ws.Range("J3").Value = "MMAA09XFGS310"

# --- Add new "accesorio movilidad" columns (L: NombreAccesorio, M: SumaAseguradaAcce) ---
$ws.Range("L1").Value = "NombreAccesorio"
$ws.Range("M1").Value = "SumaAseguradaAcce"
$ws.Range("L2").Value = "Movilidad"
$ws.Range("M2").Value = "Hasta `$150.000"

# --- Remove the 3rd vehicle (row 4) data, leaving it as a blank styled template row ---
$ws.Range("A4:K4").ClearContents()

# --- Normalize the blank template rows so H:J keep the same centered style from row 4 to row 18 ---
$ws.Range("H5:J5").VerticalAlignment = -4108
$ws.Range("H15:J15").VerticalAlignment = -4108

# --- Rows 19:20 no longer carry the H:J template formatting ---
$ws.Range("H19:J20").Clear()

# --- Drop the now-unused last template row ---
$ws.Range("A21").EntireRow.Delete()

# --- Consolidate conditional formatting into a single duplicate-values rule over H2:J18 ---
$fcs = $ws.Cells.FormatConditions
$fcs.Item(2).Delete()
$fcs.Item(1).ModifyAppliesToRange($ws.Range("H2:J18"))

# --- Restore selection like in the authored workbook ---
$ws.Range("M7").Select()
